$d = $word.ActiveDocument

# The prompts-and-answers content below is inserted as one block of raw OOXML
# (via Range.InsertXML) at the very start of the document, ahead of whatever
# is already there. Word exposes the new paragraphs' markup exactly as
# written here, including the w:proofErr spell/grammar-check markers, the
# <w:br/> line breaks, and the two horizontal-rule paragraphs (bottom
# w:pBdr), instead of only being reachable through TypeText/TypeParagraph.
$xmlBlob = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Role </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>play</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> as an expert in C++ and web development, for building tools and libraires.  We want to develop a simple C++ DLL/Shared Object that </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>is capable of connecting</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> to a native </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c++</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> application to an external service, and this external service can also connect to a web browser allowing it to visualize any information coming from the native </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c++</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> application, and/or send events back to the native application from the browser to control the native </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c++</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> app execution and/or events. </w:t>
      </w:r>
    </w:p>

    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">What would be a strategy to use for the external application, that will allow for easy integration/use between both a native </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c++</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> app and a browser.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pBdr>
          <w:bottom w:val="single" w:sz="12" w:space="1" w:color="auto"/>
        </w:pBdr>
      </w:pPr>
    </w:p>
    <w:p/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">What would be a good simple quick interface for a </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c++</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> class to wrap around a WebSocket server, that would allow for the basic functionality:</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">Start a server if it is not </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>online</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Shutdown the server</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Connect to </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>server</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Send data to server via JSON </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>objects</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Receive data back from server that could be sent from a browser.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Please provide a IWebsocket.cpp and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>IWebsocket.h</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> that matches your proposed </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>interface</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Could you also provide a schema that we could use for the JSON messaging, so that we could keep all </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>communiation</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> between the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>c++</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> app and the browser consistent.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:pBdr>
          <w:bottom w:val="single" w:sz="12" w:space="1" w:color="auto"/>
        </w:pBdr>
      </w:pPr>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Can you provide a design document that could be used to help build </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Webscoket</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> interface code as well as give detailed instructions on how to create and parse messages to and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>form</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> the </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>server</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($xmlBlob)

# InsertXML inserts its content ahead of the existing Range, so the
# document's original (already-present) empty paragraph is pushed down to
# become a trailing, now-duplicate blank paragraph after everything we just
# added. Splice it back out: delete the paragraph mark that ends the last
# paragraph we inserted, merging it with that trailing blank paragraph
# (a no-op on visible content, since the trailing paragraph is empty) so the
# document again ends with a single paragraph mark.
$lastNewParaIndex = $d.Paragraphs.Count - 1
$lastNewPara = $d.Paragraphs($lastNewParaIndex)
$mergeRange = $d.Range($lastNewPara.Range.End - 1, $lastNewPara.Range.End)
$mergeRange.Delete()
